$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48; existing rows 48:61 shift down to 49:62,
# inheriting the row's own formatting (Excel default "format as above").
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new data record.
$ws.Cells.Item(48, 1).Value2  = 11
$ws.Cells.Item(48, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(48, 3).Value2  = "Bíobío"
$ws.Cells.Item(48, 4).Value2  = 44876
$ws.Cells.Item(48, 5).Value2  = 8
$ws.Cells.Item(48, 6).Value2  = 100112031
$ws.Cells.Item(48, 7).Value2  = "Poroto verde"
$ws.Cells.Item(48, 8).Value2  = "Magnum"
$ws.Cells.Item(48, 9).Value2  = "Primera"
$ws.Cells.Item(48, 10).Value2 = 70
$ws.Cells.Item(48, 11).Value2 = 32000
$ws.Cells.Item(48, 12).Value2 = 34000
$ws.Cells.Item(48, 13).Value2 = 32857
$ws.Cells.Item(48, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(48, 15).Value2 = "Perú"
$ws.Cells.Item(48, 16).Value2 = 1314
$ws.Cells.Item(48, 17).Value2 = 25
$ws.Cells.Item(48, 18).Value2 = "Hortaliza"
